$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Change 1: "Python 3.7 or higher" -> "Python 3.8 or higher"
# -----------------------------------------------------------------------
$d.Content.Find.Execute("Python 3.7 or higher", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Python 3.8 or higher", 2) | Out-Null

# -----------------------------------------------------------------------
# Helper: locate paragraphs scoped to the "Install Python Dependencies"
# section (between its own heading and the next "5. Model Training"
# heading) so text matching does not accidentally hit other similar
# text elsewhere in the document.
# -----------------------------------------------------------------------
function Get-SectionBounds {
    $h1 = $null
    $h2 = $null
    foreach ($pp in $d.Paragraphs) {
        if ($pp.Range.Text -like "*Install Python Dependencies*") { $h1 = $pp }
        if ($pp.Range.Text -like "*5. Model Training*") { $h2 = $pp }
    }
    return @($h1.Range.End, $h2.Range.Start)
}

function Get-ScopedParagraph([string]$pattern) {
    $bounds = Get-SectionBounds
    $lo = $bounds[0]
    $hi = $bounds[1]
    $result = $null
    foreach ($pp in $d.Paragraphs) {
        if ($pp.Range.Start -ge $lo -and $pp.Range.End -le $hi -and $pp.Range.Text -like $pattern) {
            $result = $pp
        }
    }
    return $result
}

# -----------------------------------------------------------------------
# Change 2 + 3: rewrite the whole "It is recommended..." paragraph plus
# the venv/activate/pip SourceCode paragraph into the new 9-paragraph
# sequence. We select a Range starting at the (plain-styled) "It is
# recommended..." paragraph through to the end of the SourceCode
# paragraph, and assign its .Text in one shot: every newly produced
# paragraph inherits the plain (no rPr / FirstParagraph pStyle) identity
# of that very first run, which gives us a clean slate to re-apply the
# correct paragraph styles and run styles on top of.
# -----------------------------------------------------------------------
$pFirst = Get-ScopedParagraph("*It is recommended*")
$pLast  = Get-ScopedParagraph("*requirements.txt*")
$rng = $d.Range($pFirst.Range.Start, $pLast.Range.End)

$full = "It is recommended to use a virtual environment with Python 3.8+:" + "`r" + `
        "python -m venv venv" + "`r" + `
        "If the above command does not work, remove the created venv folder and run:" + "`r" + `
        "python -m venv venv --symlinks" + "`r" + `
        "On Windows, activate with:" + "`r" + `
        "venv\Scripts\activate" + "`r" + `
        "On macOS/Linux, activate with:" + "`r" + `
        "source venv/bin/activate" + "`r" + `
        "Then install dependencies:" + "`r" + `
        "pip install -r requirements.txt"

$rng.Text = $full

# -----------------------------------------------------------------------
# Re-acquire the 9 new paragraphs (in document order) and fix up their
# paragraph styles.
# -----------------------------------------------------------------------
$newParas = @()
foreach ($pp in $d.Paragraphs) {
    if ($pp.Range.Start -ge $rng.Start -and $pp.Range.End -le ($rng.Start + $full.Length + 1)) {
        $newParas += $pp
    }
}

$pStyles = @("FirstParagraph","SourceCode","FirstParagraph","SourceCode","FirstParagraph","SourceCode","FirstParagraph","SourceCode","FirstParagraph","SourceCode")
# index 0 is the "It is recommended..." paragraph - keep FirstParagraph (no-op) then the rest follow.
for ($i = 0; $i -lt $newParas.Count; $i++) {
    $newParas[$i].Style = $pStyles[$i]
}

# -----------------------------------------------------------------------
# Apply character (run) styles to the specific sub-strings that need
# them, leaving everything else with no run style (plain), matching the
# target structure.
# -----------------------------------------------------------------------
function Set-SubStyle($paragraph, [string]$needle, [string]$styleName) {
    $t = $paragraph.Range.Text
    $found = $t.IndexOf($needle)
    if ($found -lt 0) {
        Write-Host "WARNING: substring not found: [$needle] in [$t]"
        return
    }
    $s = $paragraph.Range.Start + $found
    $e = $s + $needle.Length
    $d.Range($s, $e).Style = $styleName
}

# Paragraph: "python -m venv venv"  (ExtensionTok/NormalTok/AttributeTok/NormalTok)
$pB = $newParas[1]
Set-SubStyle $pB "python" "ExtensionTok"
Set-SubStyle $pB " " "NormalTok"
Set-SubStyle $pB "-m" "AttributeTok"
Set-SubStyle $pB " venv venv" "NormalTok"

# Paragraph: "If the above command does not work, remove the created venv folder and run:"
$pC = $newParas[2]
Set-SubStyle $pC "venv" "VerbatimChar"

# Paragraph: "python -m venv venv --symlinks"  (ExtensionTok/NormalTok/AttributeTok/NormalTok/AttributeTok)
$pD = $newParas[3]
Set-SubStyle $pD "python" "ExtensionTok"
Set-SubStyle $pD " " "NormalTok"
Set-SubStyle $pD "-m" "AttributeTok"
Set-SubStyle $pD " venv venv " "NormalTok"
Set-SubStyle $pD "--symlinks" "AttributeTok"

# Paragraph: "venv\Scripts\activate" (ExtensionTok)
$pF = $newParas[5]
Set-SubStyle $pF "venv\Scripts\activate" "ExtensionTok"

# Paragraph: "source venv/bin/activate" (BuiltInTok/NormalTok)
$pH = $newParas[7]
Set-SubStyle $pH "source" "BuiltInTok"
Set-SubStyle $pH " venv/bin/activate" "NormalTok"

# Paragraph: "pip install -r requirements.txt" (ExtensionTok/NormalTok/AttributeTok/NormalTok)
$pJ = $newParas[9]
Set-SubStyle $pJ "pip" "ExtensionTok"
Set-SubStyle $pJ " install " "NormalTok"
Set-SubStyle $pJ "-r" "AttributeTok"
Set-SubStyle $pJ " requirements.txt" "NormalTok"

Write-Host "Done."
